$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Jabłka / 123) was a duplicate of row 5 - delete it entirely,
# shifting rows 7-8 up to become rows 6-7.
$ws.Rows("6").Delete()

# Update the active selection to match the saved view state.
$ws.Range("D6").Select()
